# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 43

# Header cells in columns AD, AE, AF (mirrors header style of existing columns).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy formatting from an existing header cell so the new headers match the
# bold/bordered/centered look of the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the team record values for every data row (2 through lastRow).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 82   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 80   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
